# Actualización al 21 de mayo
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Ingreso" (sheet1): 7 new contribution rows (425-431)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Ingreso")

$ingreso = @(
    @{ row = 425; member = "Carlos";  amount = 50 },
    @{ row = 426; member = "Joel";    amount = 100 },
    @{ row = 427; member = "Orlando"; amount = 100 },
    @{ row = 428; member = "Randy";   amount = 100 },
    @{ row = 429; member = "Julio";   amount = 100 },
    @{ row = 430; member = "Alfredo"; amount = 100 },
    @{ row = 431; member = "Anuel";   amount = 100 }
)

foreach ($item in $ingreso) {
    $r = $item.row
    $ws1.Range("A$r").Value = 45067
    $ws1.Range("B$r").Value = $item.member
    $ws1.Range("C$r").Value = $item.amount
    $ws1.Range("C$r").Style = "Normal"
    $ws1.Range("D$r").Value = "Aporte"
}

[void]$ws1.Range("C428").Select()

# ---------------------------------------------------------------
# Sheet "Cuentas por cobrar" (sheet3): 1 new receivable row (3)
# (populated before "Gastos" below so new shared strings line up
#  with the order they were typed in the original editing session)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Cuentas por cobrar")

[void]$ws3.Range("A2").Copy()
[void]$ws3.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("A3").Value = 45067
$ws3.Range("B3").Value = "Kibelo"
$ws3.Range("C3").Value = "Tecnica"
$ws3.Range("D3").Value = 100
$ws3.Range("F3").Value = "Cogió un pique y picó la pelota muy duro"

[void]$ws3.Range("A4").Select()

# ---------------------------------------------------------------
# Sheet "Gastos" (sheet2): 3 new expense rows (41-43)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gastos")

$ws2.Range("A41").Value = 45067
$ws2.Range("B41").Value = "Agua y hielo"
$ws2.Range("C41").Value = 140

$ws2.Range("A42").Value = 45067
$ws2.Range("B42").Value = "Empanadas"
$ws2.Range("C42").Value = 900

$ws2.Range("A43").Value = 45067
$ws2.Range("B43").Value = "Aporte mono (pintar cancha)"
$ws2.Range("C43").Value = 1000

[void]$ws2.Range("B42").Select()

# Restore "Ingreso" as the active/visible tab (matches the saved
# workbook state) without disturbing its own remembered selection.
[void]$ws1.Activate()

Write-Host "Actualizacion al 21 de mayo aplicada"
